$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Question 6 text currently lives in several runs (the grammar
#    checker split it around `@app.route("/")`). Collapse it back to a
#    single run with a single <w:t>. Find.Execute's replacement text
#    goes through Word's smart-quote AutoFormat, so the straight quotes
#    around "/" would get curled if we put them directly in the
#    replacement string. Work around it by replacing with a
#    quote-free placeholder first, then swapping the placeholder for a
#    literal straight quote via a plain Range.Text assignment (which,
#    unlike Find.Execute's replacement, is not auto-corrected).
# ---------------------------------------------------------------------
$q6Old = '6. In Flask, what does the `@app.route("/")` decorator do?'
$q6Placeholder = '6. In Flask, what does the `@app.route(QUOTEMARKER/QUOTEMARKER)` decorator do?'
$found6 = $d.Content.Find.Execute($q6Old, $false, $false, $false, $false, $false, $true, 1, $false, $q6Placeholder, 2)
if (-not $found6) {
    Write-Host "WARNING: question 6 text not found"
}

$quoteRng = $d.Content
$quoteRng.Find.ClearFormatting()
$quoteRng.Find.Text = "QUOTEMARKER"
$quoteRng.Find.Forward = $true
$quoteRng.Find.Wrap = 1
while ($quoteRng.Find.Execute()) {
    $quoteRng.Text = [string][char]0x22
    $quoteRng.Collapse(0)
}

# ---------------------------------------------------------------------
# 2. Question 10 text is also split across runs (around "recommend
#    against"); it has no special characters so a straight
#    Find/Replace is safe.
# ---------------------------------------------------------------------
$q10Old = '10. Outline the pros and cons of using Flask for a new project. Under what circumstances would you recommend using Flask, and when might you recommend against it?'
$found10 = $d.Content.Find.Execute($q10Old, $false, $false, $false, $false, $false, $true, 1, $false, $q10Old, 2)
if (-not $found10) {
    Write-Host "WARNING: question 10 text not found"
}

# ---------------------------------------------------------------------
# 3. Trim five of the blank paragraphs between question 10 and the
#    "Answer Key" heading.
# ---------------------------------------------------------------------
$q10Para = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq $q10Old) {
        $q10Para = $p
        break
    }
}
if ($null -eq $q10Para) {
    Write-Host "WARNING: could not locate question 10 paragraph"
} else {
    $removeCount = 5
    for ($i = 0; $i -lt $removeCount; $i++) {
        $blank = $q10Para.Next()
        $blank.Range.Delete()
    }
}

# ---------------------------------------------------------------------
# 4. Stamp the "Answer Key" heading with a lastRenderedPageBreak, as it
#    now starts a fresh page.
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Answer Key") {
        $r = $p.Range
        $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="36CA5882" w14:textId="65D1BE48" w:rsidR="008C65A8" w:rsidRDefault="00000000"><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Answer Key</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
        $r.InsertXML($xml)
        break
    }
}
